$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1116
$ws1.Range("F8").Value = 398
$ws1.Range("F11").Value = 511
$ws1.Range("F12").Value = 539
$ws1.Range("F14").Value = 12742
$ws1.Range("F16").Value = 5245
$ws1.Range("F17").Value = 5528

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 98

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1116
$ws4.Range("F9").Value = 398
$ws4.Range("F12").Value = 511
$ws4.Range("F13").Value = 539
$ws4.Range("F15").Value = 12742
$ws4.Range("F16").Value = 98
$ws4.Range("F19").Value = 5245
$ws4.Range("F20").Value = 5528
